$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Split paragraphs that need to become their own paragraph.
# We use Find/Replace with a literal paragraph mark (^p) in the replacement
# text, which is the standard Word trick for splitting a paragraph in two.
# Doing this first (before touching formatting) keeps paragraph indices easy
# to reason about afterwards.
# ---------------------------------------------------------------------------

# Split "VISTO: " away from the rest of paragraph 3.
$rng = $d.Content
$rng.Find.Execute("VISTO: La", $true, $false, $false, $false, $false, $true, 1, $false, "VISTO: ^p La", 2)

# Split "CONSIDERANDO: " away from the rest (now paragraph 4).
$rng = $d.Content
$rng.Find.Execute("CONSIDERANDO: Que", $true, $false, $false, $false, $false, $true, 1, $false, "CONSIDERANDO: ^p Que", 2)

# Collapse the long run of spaces before "(tres)" down to a single space.
$rng = $d.Content
$rng.Find.Execute("                    (tres)", $true, $false, $false, $false, $false, $true, 1, $false, " (tres)", 2)

# Drop the leading "POR " before "EL CONCEJO DELIBERANTE..."
$rng = $d.Content
$rng.Find.Execute("POR EL CONCEJO DELIBERANTE", $true, $false, $false, $false, $false, $true, 1, $false, "EL CONCEJO DELIBERANTE", 2)

# Remove the stray leading-space run at the start of the "Que en sus..." and
# "Que se debe..." paragraphs.
$rng = $d.Content
$rng.Find.Execute(" Que en sus Artículos", $true, $false, $false, $false, $false, $true, 1, $false, "Que en sus Artículos", 2)
$rng = $d.Content
$rng.Find.Execute(" Que se debe corregir", $true, $false, $false, $false, $false, $true, 1, $false, "Que se debe corregir", 2)

# "ARTICULO PRIMERO/SEGUNDO" headings: the leading space that used to sit
# before the heading now sits right after the colon instead.
$rng = $d.Content
$rng.Find.Execute(" ARTICULO PRIMERO: MODIFICA", $true, $false, $false, $false, $false, $true, 1, $false, "ARTICULO PRIMERO: MODIFICA", 2)
$rng = $d.Content
$rng.Find.Execute(" ARTICULO SEGUNDO: COMUNIQUESE", $true, $false, $false, $false, $false, $true, 1, $false, "ARTICULO SEGUNDO: COMUNIQUESE", 2)

Write-Output ("ParaCount after text edits=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    Write-Output ($i.ToString() + " len=" + $p.Range.Text.Length)
}

# ---------------------------------------------------------------------------
# Step 2: paragraph-level formatting
# ---------------------------------------------------------------------------

# P1: "Yerba Buena, 19 de Junio de 2014"
$p = $d.Paragraphs(1)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 12        # 240 twips = 12pt

# P2: "ORDENANZA Nº 1983"
$p = $d.Paragraphs(2)
$p.Format.KeepWithNext = $true
$p.Format.SpaceBefore = 12       # 240 twips = 12pt
$p.Format.SpaceAfter = 18        # 360 twips = 18pt
$p.Range.Font.Bold = $true

# P3: "VISTO: "
$p = $d.Paragraphs(3)
$p.Format.KeepWithNext = $true
$p.Format.SpaceBefore = 12       # 240 twips = 12pt
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # wdAlignParagraphLeft
$p.Range.Font.Bold = $true

# P4: " La Ordenanza Nº 1962..."
$p = $d.Paragraphs(4)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # wdAlignParagraphLeft

# P5: "CONSIDERANDO: "
$p = $d.Paragraphs(5)
$p.Format.KeepWithNext = $true
$p.Format.SpaceBefore = 12       # 240 twips = 12pt
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # wdAlignParagraphJustify -> none in target (removed jc=both)
$p.Range.Font.Bold = $true

# P6: " Que a través..."
$p = $d.Paragraphs(6)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # clears jc="both"

# P7: "Que en sus Artículos..."
$p = $d.Paragraphs(7)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # clears jc="both"

# P8: "Que se debe corregir..."
$p = $d.Paragraphs(8)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # clears jc="both"

# P9: "EL CONCEJO DELIBERANTE..."
$p = $d.Paragraphs(9)
$p.Format.KeepWithNext = $true
$p.Format.SpaceBefore = 18       # 360 twips = 18pt
$p.Format.SpaceAfter = 18        # 360 twips = 18pt
$p.Format.LeftIndent = 99.2      # 1984 twips = 99.2pt (1/20 pt per twip)
$p.Format.RightIndent = 99.2
$p.Range.Font.Bold = $true

# P10: "ARTICULO PRIMERO: MODIFICANSE..."
$p = $d.Paragraphs(10)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 6         # 120 twips = 6pt

# P11: "ARTICULO SEGUNDO: COMUNIQUESE..."
$p = $d.Paragraphs(11)
$p.Format.KeepWithNext = $true
$p.Format.SpaceAfter = 6         # 120 twips = 6pt
$p.Format.Alignment = 0          # clears jc="both"

# ---------------------------------------------------------------------------
# Step 3: underline the "ARTICULO PRIMERO:" / "ARTICULO SEGUNDO:" headings
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("ARTICULO PRIMERO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 1

$rng = $d.Content
$rng.Find.Execute("ARTICULO SEGUNDO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 1

Write-Output ("ParaCount after splits=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    Write-Output ($i.ToString() + " len=" + $p.Range.Text.Length)
}
